$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New simulated run values (log write mode update) for rows 2-11,
# columns C (run_time), E (max_er), and G:Y (iter 1..19).
$updates = @{
    "C2" = 0.4758577346801758
    "E2" = 1736.777416304632
    "G2" = 0.1393610345579782
    "H2" = 0.1158110994474369
    "I2" = 0.09834253961323451
    "J2" = 0.08196577634368189
    "K2" = 0.06779934550973617
    "L2" = 0.05811601412699177
    "M2" = 0.05526666972479417
    "N2" = 0.05015347686396227
    "O2" = 0.04528585763464777
    "P2" = 0.0429545655320468
    "Q2" = 0.03989275682576304
    "R2" = 0.03849160907873388
    "S2" = 0.03763699137270304
    "T2" = 0.03585776277715585
    "U2" = 0.03552068291220589
    "V2" = 0.035082147586712
    "W2" = 0.03434355642023491
    "X2" = 0.03414453758935892
    "Y2" = 0.03385531025934956
    "C3" = 0.3918361663818359
    "E3" = 1737.972495023905
    "G3" = 0.1439895923815842
    "H3" = 0.1397108826941851
    "I3" = 0.1077081401097672
    "J3" = 0.07987496355142751
    "K3" = 0.06424011989373378
    "L3" = 0.05700209398889326
    "M3" = 0.05122717666982481
    "N3" = 0.04627806204020907
    "O3" = 0.04416243213229053
    "P3" = 0.0409257991822928
    "Q3" = 0.03966439653062103
    "R3" = 0.03764268546429202
    "S3" = 0.03695453607751747
    "T3" = 0.0355225952297849
    "U3" = 0.0348296126002086
    "V3" = 0.03447756685162036
    "W3" = 0.034174653638549
    "X3" = 0.03402216285593923
    "Y3" = 0.03387860614081686
    "C4" = 0.3906011581420898
    "E4" = 1961.055063488344
    "G4" = 0.1382372116903636
    "H4" = 0.1328849675760826
    "I4" = 0.1097573795898677
    "J4" = 0.09142047352277943
    "K4" = 0.07579808568973788
    "L4" = 0.06447478869502804
    "M4" = 0.05832333558387
    "N4" = 0.05307123559270142
    "O4" = 0.04876648318049965
    "P4" = 0.04632406014169593
    "Q4" = 0.04465072280025808
    "R4" = 0.04311683431140344
    "S4" = 0.04170557859060321
    "T4" = 0.04063424123876128
    "U4" = 0.03982963897308833
    "V4" = 0.03909886063622861
    "W4" = 0.03855954911622538
    "X4" = 0.03835832129541582
    "Y4" = 0.03822719422004568
    "C5" = 0.3906214237213135
    "E5" = 1764.840653124647
    "G5" = 0.1471089493055333
    "H5" = 0.1368444988783983
    "I5" = 0.1007854893198985
    "J5" = 0.07903864948513337
    "K5" = 0.06200886202792739
    "L5" = 0.0536775941027101
    "M5" = 0.04876156976287127
    "N5" = 0.04627233603268054
    "O5" = 0.04251211287139038
    "P5" = 0.0393475702001328
    "Q5" = 0.03813028938420667
    "R5" = 0.03725430146838997
    "S5" = 0.03616415374188189
    "T5" = 0.03574724009609337
    "U5" = 0.03534281106333326
    "V5" = 0.03510461990035327
    "W5" = 0.03475354837941966
    "X5" = 0.03450498263028193
    "Y5" = 0.03440235191276115
    "C6" = 0.4531536102294922
    "E6" = 1789.807811070305
    "G6" = 0.1416509177669774
    "H6" = 0.136908844838888
    "I6" = 0.09072137281266347
    "J6" = 0.07278515194162195
    "K6" = 0.06604491353509478
    "L6" = 0.0568828950153947
    "M6" = 0.05235470777852642
    "N6" = 0.04742469786317166
    "O6" = 0.04458003458389231
    "P6" = 0.04151551324950774
    "Q6" = 0.04063677669463553
    "R6" = 0.03870711664312331
    "S6" = 0.03770277160247686
    "T6" = 0.0370282817226416
    "U6" = 0.03632759970741266
    "V6" = 0.03586669630806247
    "W6" = 0.03528501864275206
    "X6" = 0.0349978976050932
    "Y6" = 0.03488904115146792
    "C7" = 0.3905997276306152
    "E7" = 1799.606201773744
    "G7" = 0.14803351255148
    "H7" = 0.1414548836233981
    "I7" = 0.1315209565581311
    "J7" = 0.1057899072818349
    "K7" = 0.09050900587488052
    "L7" = 0.07092943739956918
    "M7" = 0.05998072753450132
    "N7" = 0.05360051798347471
    "O7" = 0.0490645215358218
    "P7" = 0.04601582668630021
    "Q7" = 0.04249468619943188
    "R7" = 0.04064952487242499
    "S7" = 0.03911461217487257
    "T7" = 0.03742011255251643
    "U7" = 0.03633264027287295
    "V7" = 0.03596542506226681
    "W7" = 0.03559985917040058
    "X7" = 0.03529331799435547
    "Y7" = 0.03508004291956615
    "C8" = 0.3750374317169189
    "E8" = 1892.063820495421
    "G8" = 0.1488505658856033
    "H8" = 0.1412161107205949
    "I8" = 0.1231876455875655
    "J8" = 0.1004588351315294
    "K8" = 0.07867778647224154
    "L8" = 0.06891218197116693
    "M8" = 0.06136601031642046
    "N8" = 0.05197450967238423
    "O8" = 0.04910462698947399
    "P8" = 0.04560332177923785
    "Q8" = 0.04404630417163034
    "R8" = 0.04196257649508228
    "S8" = 0.04039021101141394
    "T8" = 0.03900301454862413
    "U8" = 0.03841251464219408
    "V8" = 0.03784669428657016
    "W8" = 0.03734082969009179
    "X8" = 0.03710840470426585
    "Y8" = 0.03688233568217195
    "C9" = 0.3749594688415527
    "E9" = 1995.042158747146
    "G9" = 0.1470971487920564
    "H9" = 0.1370668248027466
    "I9" = 0.11845601410399
    "J9" = 0.09309888386956357
    "K9" = 0.07664180059211662
    "L9" = 0.06757306216111846
    "M9" = 0.05843204675792778
    "N9" = 0.05400920691651272
    "O9" = 0.05054557133252551
    "P9" = 0.04675818708243432
    "Q9" = 0.04440036811163021
    "R9" = 0.04365427967400145
    "S9" = 0.04221208736245347
    "T9" = 0.04143570754934179
    "U9" = 0.04070529467988147
    "V9" = 0.040081591383836
    "W9" = 0.03941763700812664
    "X9" = 0.03914131974370871
    "Y9" = 0.03888971069682546
    "C10" = 0.3750274181365967
    "E10" = 1898.312054279198
    "G10" = 0.1447155221672362
    "H10" = 0.1404046849725465
    "I10" = 0.1080707280037227
    "J10" = 0.08540842642878879
    "K10" = 0.07202244574391832
    "L10" = 0.06150003956882644
    "M10" = 0.05542082687165746
    "N10" = 0.05151144611095325
    "O10" = 0.04800215036336795
    "P10" = 0.04455389486524126
    "Q10" = 0.0415243019114908
    "R10" = 0.04114694330997944
    "S10" = 0.04032926253779493
    "T10" = 0.03939823907699906
    "U10" = 0.03829791466883761
    "V10" = 0.03783809452427127
    "W10" = 0.03743619056428547
    "X10" = 0.03721766950425171
    "Y10" = 0.03700413361168027
    "C11" = 0.3593757152557373
    "E11" = 1766.466700830066
    "G11" = 0.1497792556575515
    "H11" = 0.1346625434146599
    "I11" = 0.1030318997850675
    "J11" = 0.08144779797735495
    "K11" = 0.06448134721077696
    "L11" = 0.0545874932276224
    "M11" = 0.05068459745888523
    "N11" = 0.04687002183324998
    "O11" = 0.043582714077805
    "P11" = 0.0408557303592063
    "Q11" = 0.0379009505497753
    "R11" = 0.03762978249268869
    "S11" = 0.03682752400838765
    "T11" = 0.036227459899469
    "U11" = 0.0356611144404331
    "V11" = 0.03517475713489578
    "W11" = 0.03479589312664136
    "X11" = 0.03443404874912408
    "Y11" = 0.03443404874912408
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"
